# Convert the two M2Doc fields (" m:userdoc OrderedSet{'zone1', 'zone2'} "
# and " m:enduserdoc ") from real Word fields (fldChar begin/instrText/
# fldChar end) into plain literal text runs wrapped in "{" / "}", as
# produced by the TokenIteratorFieldRewriterSplit parser.

$d = $word.ActiveDocument

# --- Field 1: " m:userdoc OrderedSet{'zone1', 'zone2'} " -> plain text ---
$f1 = $d.Fields(1)
$p1 = $f1.Code.Paragraphs(1)
$start1 = $p1.Range.Start
$f1.Delete()

$ins1 = $d.Range($start1, $start1)
$ins1.InsertBefore("{m:userdoc OrderedSet{'zone1', 'zone2'}}")

# --- Field 2: " m:enduserdoc " -> plain text ---
$f2 = $d.Fields(1)
$p2 = $f2.Code.Paragraphs(1)
$start2 = $p2.Range.Start
$f2.Delete()

$ins2 = $d.Range($start2, $start2)
$ins2.InsertBefore("{m:enduserdoc}")
